# "Generate Report for Archive"
#
# The localization status of the single tracked file moved from
# "Ready for handoff" to "In Translation". That shared-string value is
# surfaced on all three worksheets:
#   - Overview : columns "zh-cn" (E2) and "de-de" (F2)
#   - zh-cn    : "Status" column (C2)
#   - de-de    : "Status" column (C2)
#
# Shortening that status text is also what made Excel shrink the
# (now too-wide) "zh-cn"/"de-de"/"Status" columns on the next save, so we
# nudge those column widths down to match as closely as this host's
# ColumnWidth quantization allows.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) and de-de (col F) status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (col C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (col C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
